# Automatic update of files.
# Rotates taxon/observation data among rows 13/14/15 and 22/23/24.
# (Row identities stay put; the species-record + location-metric columns
# that belong to each observation move to a different physical row.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13 ---
$ws.Range("A13").Value = 130864511
$ws.Range("B13").Value = 83223
$ws.Range("E13").Value = 6440
$ws.Range("F13").Value = "Vitgrynig nållav"
$ws.Range("G13").Value = "Chaenotheca subroscida"
$ws.Range("H13").Value = "(Eitner) Zahlbr."
$ws.Range("Q13").Value = 445985
$ws.Range("R13").Value = 7031157
$ws.Range("S13").Value = 3
$ws.Range("Z13").Value = "10:40"
$ws.Range("AB13").Value = "10:40"

# --- Row 14 ---
$ws.Range("A14").Value = 130864514
$ws.Range("Q14").Value = 445944
$ws.Range("R14").Value = 7031101
$ws.Range("S14").Value = 4
$ws.Range("Z14").Value = "11:03"
$ws.Range("AB14").Value = "11:03"

# --- Row 15 ---
$ws.Range("A15").Value = 130864513
$ws.Range("B15").Value = 78255
$ws.Range("E15").Value = 228579
$ws.Range("F15").Value = "Liten svartspik"
$ws.Range("G15").Value = "Chaenothecopsis nana"
$ws.Range("H15").Value = "Tibell"
$ws.Range("Q15").Value = 445962
$ws.Range("R15").Value = 7031114
$ws.Range("Z15").Value = "10:51"
$ws.Range("AB15").Value = "10:51"

# --- Row 22 ---
$ws.Range("A22").Value = 130864510
$ws.Range("B22").Value = 78255
$ws.Range("E22").Value = 228579
$ws.Range("F22").Value = "Liten svartspik"
$ws.Range("G22").Value = "Chaenothecopsis nana"
$ws.Range("H22").Value = "Tibell"
$ws.Range("Q22").Value = 445990
$ws.Range("R22").Value = 7031220
$ws.Range("S22").Value = 6
$ws.Range("Z22").Value = "10:33"
$ws.Range("AB22").Value = "10:33"
$ws.Range("AC22").Value = "På bark på stam av levande gammal gran i gammal granskog"

# --- Row 23 ---
$ws.Range("A23").Value = 130864515
$ws.Range("B23").Value = 83223
$ws.Range("D23").Value = "NT"
$ws.Range("E23").Value = 6440
$ws.Range("F23").Value = "Vitgrynig nållav"
$ws.Range("G23").Value = "Chaenotheca subroscida"
$ws.Range("H23").Value = "(Eitner) Zahlbr."
$ws.Range("Q23").Value = 445932
$ws.Range("R23").Value = 7031103
$ws.Range("S23").Value = 3
$ws.Range("Z23").Value = "11:08"
$ws.Range("AB23").Value = "11:08"
$ws.Range("AC23").Value = "På bark på stam av levande gammal gran"

# --- Row 24 ---
$ws.Range("A24").Value = 130864521
$ws.Range("B24").Value = 91771
$ws.Range("D24").Value = "LC"
$ws.Range("E24").Value = 5447
$ws.Range("F24").Value = "Vedticka"
$ws.Range("G24").Value = "Fuscoporia viticola"
$ws.Range("H24").Value = "(Schwein.) Murrill"
$ws.Range("Q24").Value = 446069
$ws.Range("R24").Value = 7030939
$ws.Range("S24").Value = 8
$ws.Range("Z24").Value = "13:20"
$ws.Range("AB24").Value = "13:20"
$ws.Range("AC24").Value = "På död klen gran i gammal granskog"

# K24/N24/AF24 become blank, K23/N23/AF23 are now present as blanks.
$ws.Range("K24").ClearContents()
$ws.Range("N24").ClearContents()
$ws.Range("AF24").ClearContents()
